$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeiterfassung")

# Append new time-tracking rows (19-22)
$ws.Range("A19").Value = 45597
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = "Coding"
$ws.Range("D19").Value = "Objektparser programmieren"

$ws.Range("A20").Value = 45598
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = "Coding"
$ws.Range("D20").Value = "Objektparser fertigstellen & Materialparser programmieren"

$ws.Range("A21").Value = 45599
$ws.Range("B21").Value = 2.5
$ws.Range("C21").Value = "Coding"
$ws.Range("D21").Value = "Materialparser fertigstellen"

$ws.Range("A22").Value = 45600
$ws.Range("B22").Value = 1.5
$ws.Range("C22").Value = "Online-Meeting"
$ws.Range("D22").Value = "Weekly-Summup-03 Meeting"

# Match the formatting of the preceding rows (date format in column A, left-aligned in column B)
$ws.Range("A18").Copy()
$ws.Range("A19:A22").PasteSpecial(-4122)
$ws.Range("B18").Copy()
$ws.Range("B19:B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to reflect the new active cell
$ws.Range("G18").Select()
